$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Definitive SMB (col Q) / dSMB 2001-2030 (col R) / dSMB 2071-2100 (col S) table,
# row 8 keeps its existing #NUM! error and is intentionally skipped.
$definitive = @(
    @(2, 1332.02174797158, 128.396594571143, 527.365239108085),
    @(3, 271.176738064172, 27.2287548961168, 90.2820686613982),
    @(4, 557.486977037644, -17.5062372651231, 60.3196963909475),
    @(5, 74.9078035724264, 2.21532775454592, 14.1874301803554),
    @(6, 379.210613351733, 57.3445139209004, 157.423922727042),
    @(7, 54.2147500658909, 5.90009467702878, 14.3674393437545),
    @(9, 137.905558514223, -0.800580838868496, -6.12978319284748),
    @(10, 308.19423633973, 1.65186434146875, 26.9844711643318),
    @(11, 142.453816661317, -3.25833330793488, 4.07545910988381),
    @(12, 75.5557164594573, 5.54666441005803, 16.6315656575833),
    @(13, 35.9912258445195, -0.0457932469920995, 10.9209203065984),
    @(14, 10.3518975583481, -1.43273131319767, 1.35985423631178),
    @(15, 947.607420535923, -30.5899985925267, 54.4490355171217),
    @(16, 1054.39414608325, -15.3534281790346, 39.3160331292686),
    @(17, 1012.33120858216, -2.93678428037042, 50.8364089339072),
    @(18, 288.911747433697, -11.5846466360928, 29.9184928620856),
    @(19, 54.2009255134342, 1.59671093444798, 20.8559894083978),
    @(20, 569.306070782311, 22.3198947270323, 156.893005449802),
    @(21, 529.378878482933, 14.4456597001585, 168.370424921899),
    @(22, 3.50147372735937, 0.0747662264202085, 0.634140848276613),
    @(23, 1172.06933673294, 21.8761862127158, -284.63297006762)
)

foreach ($entry in $definitive) {
    $r = $entry[0]
    $ws.Cells.Item($r, 17).Value = $entry[1]
    $ws.Cells.Item($r, 18).Value = $entry[2]
    $ws.Cells.Item($r, 19).Value = $entry[3]
}
